# D4_analytical_dataset_KM_ImmDis.xlsx -- "minor modifications to index"
#
# Semantic changes applied:
#  1. Metadata sheet: "Content of the dataset" (B3) text is extended with a
#     note about how the dataset is derived; row height grows to fit it.
#     Selection on that sheet moves from B9 to the full A2:B10 block.
#  2. Data Model sheet: a new row is inserted describing the new
#     "days_{ImmDis}" variable (number of days in the follow up). A couple of
#     stray notes (age_at_cohort_entry_date / start_follow_up / start_period /
#     end_period rows) are cleared out, and the Retrieved/Calculated "yes"
#     markers are straightened out onto their correct rows/columns.
#  3. The active/selected worksheet switches from "Example" to "Data Model".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "Persons that are in the cohort of {ImmDis}, with first periods of follow-up and flare (if any): this is the analytical dataset that enters the calculation of cumulative incidence at 180 days and at 365 days. This dataset is obtained from D3_followup_periods_in_cohort_{ImmDis} by selecting the first record (number_of_period_{ImmDis}) == 1"
$meta.Rows(3).RowHeight = 72

# Selection moves from B9 to A2:B10 (this also activates the Metadata sheet,
# which is fine -- a later step re-activates "Data Model" as the final tab).
$meta.Range("A2:B10").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Data Model sheet
# ---------------------------------------------------------------------
$dm = $wb.Worksheets.Item("Data Model")

# Row 3 (sex_at_instance_creation): the "yes" marker moves from the
# "Calculated" column (I) to the "Retrieved" column (H).
$dm.Range("I3").ClearContents()
$dm.Range("H3").Value = "yes"

# Row 4 (age_at_cohort_entry_date_{ImmDis}): drop the stray note, and it
# turns out to be "Calculated" (I), not left blank. Without the long note the
# row no longer needs its tall custom height.
$dm.Range("E4").ClearContents()
$dm.Range("I4").Value = "yes"
$dm.Rows(4).AutoFit() | Out-Null

# Row 5 (start_follow_up_{ImmDis}_d): "yes" moves from Calculated (I) to
# Retrieved (H); the formula note in Rule (K) is dropped, so the row shrinks
# back to its default height too.
$dm.Range("I5").ClearContents()
$dm.Range("H5").Value = "yes"
$dm.Range("K5").ClearContents()
$dm.Rows(5).AutoFit() | Out-Null

# Row 6 (start_period_{ImmDis}_d): drop the note, and add Parameters/Retrieved.
$dm.Range("E6").ClearContents()
$dm.Range("F6").Value = "ImmDis"
$dm.Range("H6").Value = "yes"

# Row 7 (end_period_{ImmDis}_d): drop the note, and add Parameters/Retrieved.
$dm.Range("E7").ClearContents()
$dm.Range("F7").Value = "ImmDis"
$dm.Range("H7").Value = "yes"

# Insert a new row for the "days_{ImmDis}" variable right before the
# cause_end_period_{ImmDis} row (pushes cause_end_period_/flare_ down by one).
$dm.Rows(8).Insert()

$dm.Range("A8").Value = "days_{ImmDis}"
$dm.Range("B8").Value = "number of days in the follow up"
$dm.Range("F8").Value = "ImmDis"
$dm.Range("I8").Value = "yes"

# cause_end_period_{ImmDis} (now row 9) gains a "Parameters" entry.
$dm.Range("F9").Value = "ImmDis"

# flare_{ImmDis} (now row 10) gains a "Parameters" entry.
$dm.Range("F10").Value = "ImmDis"

# Data Model becomes the active sheet, with the selection covering the
# populated block in the frozen bottom-right pane.
$dm.Activate() | Out-Null
$dm.Range("A2:V14").Select() | Out-Null
